$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 251-259 (only changed cells, per diff) ---
# Row 251
$ws.Range("D251").Value = 44448
$ws.Range("M251").Value = 20
$ws.Range("N251").Value = 175000
$ws.Range("O251").Value = 180000
$ws.Range("P251").Value = 177000
$ws.Range("Q251").Value = '$/bins (450 kilos)'
$ws.Range("R251").Value = 'Llay Llay'
$ws.Range("S251").Value = 393
$ws.Range("T251").Value = 450

# Row 252
$ws.Range("D252").Value = 44448
$ws.Range("L252").Value = 'Especial'
$ws.Range("M252").Value = 450
$ws.Range("N252").Value = 9000
$ws.Range("O252").Value = 9000
$ws.Range("P252").Value = 9000
$ws.Range("Q252").Value = '$/caja 18 kilos'
$ws.Range("R252").Value = 'Región Metropolitana'
$ws.Range("S252").Value = 500
$ws.Range("T252").Value = 18

# Row 253
$ws.Range("D253").Value = 44448
$ws.Range("L253").Value = 'Primera'
$ws.Range("M253").Value = 16
$ws.Range("N253").Value = 130000
$ws.Range("O253").Value = 140000
$ws.Range("P253").Value = 136250
$ws.Range("Q253").Value = '$/bins (450 kilos)'
$ws.Range("R253").Value = 'Llay Llay'
$ws.Range("S253").Value = 303
$ws.Range("T253").Value = 450

# Row 254
$ws.Range("D254").Value = 44448
$ws.Range("M254").Value = 400
$ws.Range("N254").Value = 7000
$ws.Range("O254").Value = 7000
$ws.Range("P254").Value = 7000
$ws.Range("Q254").Value = '$/caja 18 kilos'
$ws.Range("R254").Value = 'Región Metropolitana'
$ws.Range("S254").Value = 389
$ws.Range("T254").Value = 18

# Row 255
$ws.Range("D255").Value = 44448
$ws.Range("K255").Value = 'Murcott'
$ws.Range("L255").Value = 'Segunda'
$ws.Range("M255").Value = 350
$ws.Range("N255").Value = 5500
$ws.Range("O255").Value = 5500
$ws.Range("P255").Value = 5500
$ws.Range("Q255").Value = '$/caja 18 kilos'
$ws.Range("R255").Value = 'Región Metropolitana'
$ws.Range("S255").Value = 306
$ws.Range("T255").Value = 18

# Row 256
$ws.Range("D256").Value = 44167
$ws.Range("K256").Value = 'Murcott'
$ws.Range("L256").Value = 'Especial'
$ws.Range("M256").Value = 240
$ws.Range("N256").Value = 10000
$ws.Range("O256").Value = 10000
$ws.Range("P256").Value = 10000
$ws.Range("S256").Value = 1000

# Row 257
$ws.Range("D257").Value = 44167
$ws.Range("K257").Value = 'Murcott'
$ws.Range("L257").Value = 'Primera'
$ws.Range("M257").Value = 220
$ws.Range("N257").Value = 8000
$ws.Range("O257").Value = 8000
$ws.Range("P257").Value = 8000
$ws.Range("S257").Value = 800

# Row 258
$ws.Range("D258").Value = 44238
$ws.Range("K258").Value = 'Murcott'
$ws.Range("M258").Value = 60
$ws.Range("N258").Value = 12000
$ws.Range("O258").Value = 12000
$ws.Range("P258").Value = 12000
$ws.Range("S258").Value = 1200

# Row 259
$ws.Range("D259").Value = 44238
$ws.Range("K259").Value = 'Murcott'
$ws.Range("M259").Value = 90
$ws.Range("N259").Value = 10000
$ws.Range("O259").Value = 10000
$ws.Range("P259").Value = 10000
$ws.Range("S259").Value = 1000

# --- Add new rows 260-265 (full row data; row 260 pre-existing gets overwritten with shifted data, 261-265 are brand new) ---
# Row 260
$ws.Range("A260").Value = 9
$ws.Range("B260").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C260").Value = 'Metropolitana'
$ws.Range("D260").Value = 44399
$ws.Range("E260").Value = 13
$ws.Range("F260").Value = 'Fruta'
$ws.Range("G260").Value = 100102
$ws.Range("H260").Value = 'Cítricos'
$ws.Range("I260").Value = 100102004
$ws.Range("J260").Value = 'Mandarina'
$ws.Range("K260").Value = 'Clemenuless'
$ws.Range("L260").Value = 'Especial'
$ws.Range("M260").Value = 420
$ws.Range("N260").Value = 7000
$ws.Range("O260").Value = 7500
$ws.Range("P260").Value = 7262
$ws.Range("Q260").Value = '$/bandeja 10 kilos'
$ws.Range("R260").Value = 'Provincia del Elquí'
$ws.Range("S260").Value = 726
$ws.Range("T260").Value = 10

# Row 261
$ws.Range("A261").Value = 9
$ws.Range("B261").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C261").Value = 'Metropolitana'
$ws.Range("D261").Value = 44399
$ws.Range("E261").Value = 13
$ws.Range("F261").Value = 'Fruta'
$ws.Range("G261").Value = 100102
$ws.Range("H261").Value = 'Cítricos'
$ws.Range("I261").Value = 100102004
$ws.Range("J261").Value = 'Mandarina'
$ws.Range("K261").Value = 'Clemenuless'
$ws.Range("L261").Value = 'Primera'
$ws.Range("M261").Value = 500
$ws.Range("N261").Value = 5500
$ws.Range("O261").Value = 6000
$ws.Range("P261").Value = 5740
$ws.Range("Q261").Value = '$/bandeja 10 kilos'
$ws.Range("R261").Value = 'Provincia del Elquí'
$ws.Range("S261").Value = 574
$ws.Range("T261").Value = 10
$ws.Range("D261").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 262
$ws.Range("A262").Value = 9
$ws.Range("B262").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C262").Value = 'Metropolitana'
$ws.Range("D262").Value = 44399
$ws.Range("E262").Value = 13
$ws.Range("F262").Value = 'Fruta'
$ws.Range("G262").Value = 100102
$ws.Range("H262").Value = 'Cítricos'
$ws.Range("I262").Value = 100102004
$ws.Range("J262").Value = 'Mandarina'
$ws.Range("K262").Value = 'Clemenuless'
$ws.Range("L262").Value = 'Segunda'
$ws.Range("M262").Value = 550
$ws.Range("N262").Value = 3500
$ws.Range("O262").Value = 4000
$ws.Range("P262").Value = 3755
$ws.Range("Q262").Value = '$/bandeja 10 kilos'
$ws.Range("R262").Value = 'Provincia del Elquí'
$ws.Range("S262").Value = 376
$ws.Range("T262").Value = 10
$ws.Range("D262").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 263
$ws.Range("A263").Value = 9
$ws.Range("B263").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C263").Value = 'Metropolitana'
$ws.Range("D263").Value = 44400
$ws.Range("E263").Value = 13
$ws.Range("F263").Value = 'Fruta'
$ws.Range("G263").Value = 100102
$ws.Range("H263").Value = 'Cítricos'
$ws.Range("I263").Value = 100102004
$ws.Range("J263").Value = 'Mandarina'
$ws.Range("K263").Value = 'Clemenuless'
$ws.Range("L263").Value = 'Especial'
$ws.Range("M263").Value = 350
$ws.Range("N263").Value = 7500
$ws.Range("O263").Value = 7500
$ws.Range("P263").Value = 7500
$ws.Range("Q263").Value = '$/bandeja 10 kilos'
$ws.Range("R263").Value = 'Provincia del Elquí'
$ws.Range("S263").Value = 750
$ws.Range("T263").Value = 10
$ws.Range("D263").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 264
$ws.Range("A264").Value = 9
$ws.Range("B264").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C264").Value = 'Metropolitana'
$ws.Range("D264").Value = 44400
$ws.Range("E264").Value = 13
$ws.Range("F264").Value = 'Fruta'
$ws.Range("G264").Value = 100102
$ws.Range("H264").Value = 'Cítricos'
$ws.Range("I264").Value = 100102004
$ws.Range("J264").Value = 'Mandarina'
$ws.Range("K264").Value = 'Clemenuless'
$ws.Range("L264").Value = 'Primera'
$ws.Range("M264").Value = 400
$ws.Range("N264").Value = 6000
$ws.Range("O264").Value = 6000
$ws.Range("P264").Value = 6000
$ws.Range("Q264").Value = '$/bandeja 10 kilos'
$ws.Range("R264").Value = 'Provincia del Elquí'
$ws.Range("S264").Value = 600
$ws.Range("T264").Value = 10
$ws.Range("D264").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 265
$ws.Range("A265").Value = 9
$ws.Range("B265").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C265").Value = 'Metropolitana'
$ws.Range("D265").Value = 44400
$ws.Range("E265").Value = 13
$ws.Range("F265").Value = 'Fruta'
$ws.Range("G265").Value = 100102
$ws.Range("H265").Value = 'Cítricos'
$ws.Range("I265").Value = 100102004
$ws.Range("J265").Value = 'Mandarina'
$ws.Range("K265").Value = 'Clemenuless'
$ws.Range("L265").Value = 'Segunda'
$ws.Range("M265").Value = 370
$ws.Range("N265").Value = 4000
$ws.Range("O265").Value = 4000
$ws.Range("P265").Value = 4000
$ws.Range("Q265").Value = '$/bandeja 10 kilos'
$ws.Range("R265").Value = 'Provincia del Elquí'
$ws.Range("S265").Value = 400
$ws.Range("T265").Value = 10
$ws.Range("D265").NumberFormat = "YYYY-MM-DD HH:MM:SS"
